$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.296.02"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.489.68"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'321.26"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "'108.08"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").Value = "'0.522"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").Value = "'38.73"
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("D11").Value = "'0.0811"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'18.43"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "'7.12"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "2.878.42"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "2.487.19"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "47.214.16"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'12.82"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'6.61"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "0.0₃0933"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("E22").Value = "  +14.00%  "
$ws.Range("D23").Value = "'70.28"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "'245.24"
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'25.72"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.99"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'49.56"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").Value = "'20.68"
$ws.Range("E33").Value = "  +4.33%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'0.0782"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'2.92"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").Value = "'22.81"
$ws.Range("E40").Value = "  +6.35%  "
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").Value = "'117.49"
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("D44").Value = "'0.0296"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "1.984.34"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -6.17%  "
$ws.Range("D48").Value = "'9.05"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").Value = "'1.77"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("D50").Value = "'5.09"
$ws.Range("E50").Value = "  -6.74%  "
$ws.Range("D51").Value = "'56.65"
$ws.Range("E51").Value = "  +3.04%  "
